$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-12 from 2023-09-06 (45175)
# to 2023-09-14 (45183), reflecting an automatic refresh of the data.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
